# Updated cryptos list with latest price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.885.70'
$ws.Range("E2").Value = '  -1.64%  '

$ws.Range("D3").Value = '3.380.68'
$ws.Range("E3").Value = '  -0.62%  '

$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("D5").Value = '572.45'
$ws.Range("E5").Value = '  -0.70%  '

$ws.Range("D6").Value = '137.24'
$ws.Range("E6").Value = '  +0.22%  '

$ws.Range("E7").Value = '  +0.00%  '

$ws.Range("D8").Value = '3.378.63'
$ws.Range("E8").Value = '  -0.66%  '

$ws.Range("E9").Value = '  -0.98%  '

$ws.Range("D10").Value = '7.66'
$ws.Range("E10").Value = '  +2.57%  '

$ws.Range("E11").Value = '  -2.13%  '

$ws.Range("D12").Value = '0.385'
$ws.Range("E12").Value = '  -1.86%  '

$ws.Range("E13").Value = '  -0.63%  '

$ws.Range("E14").Value = '  +0.47%  '

$ws.Range("D15").Value = '26.03'
$ws.Range("E15").Value = '  +2.19%  '

$ws.Range("E16").Value = '  -2.89%  '

$ws.Range("D17").Value = '3.379.66'
$ws.Range("E17").Value = '  -0.62%  '

$ws.Range("D18").Value = '61.030.17'
$ws.Range("E18").Value = '  -1.49%  '

$ws.Range("D19").Value = '13.98'
$ws.Range("E19").Value = '  -1.76%  '

$ws.Range("B20").Value = 'Uniswap'
$ws.Range("C20").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D20").Value = '9.47'
$ws.Range("E20").Value = '  -0.19%  '

$ws.Range("B21").Value = 'Polkadot'
$ws.Range("C21").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D21").Value = '5.81'
$ws.Range("E21").Value = '  -1.05%  '

$ws.Range("D22").Value = '374.93'
$ws.Range("E22").Value = '  -3.48%  '

$ws.Range("E23").Value = '  -2.68%  '

$ws.Range("D24").Value = '3.526.82'

$ws.Range("E25").Value = '  -0.04%  '

$ws.Range("E26").Value = '  -1.60%  '

$ws.Range("D27").Value = '71.09'
$ws.Range("E27").Value = '  -0.30%  '

$ws.Range("D28").Value = '1.78'
$ws.Range("E28").Value = '  +12.18%  '

$ws.Range("D29").Value = '0.176'
$ws.Range("E29").Value = '  +9.28%  '

$ws.Range("D30").Value = '7.49'
$ws.Range("E30").Value = '  -2.82%  '

$ws.Range("D31").Value = '0.992'
$ws.Range("E31").Value = '  -0.91%  '

$ws.Range("D32").Value = '8.14'
$ws.Range("E32").Value = '  -1.90%  '

$ws.Range("E33").Value = '  -1.29%  '

$ws.Range("E34").Value = '  -0.02%  '

$ws.Range("D35").Value = '23.69'
$ws.Range("E35").Value = '  +0.68%  '

$ws.Range("E36").Value = '  -4.03%  '

$ws.Range("D37").Value = '6.88'
$ws.Range("E37").Value = '  -1.53%  '

$ws.Range("E38").Value = '  -0.53%  '

$ws.Range("D39").Value = '164.88'
$ws.Range("E39").Value = '  +1.10%  '

$ws.Range("D40").Value = '0.0762'
$ws.Range("E40").Value = '  -3.30%  '

$ws.Range("E41").Value = '  -0.08%  '

$ws.Range("E42").Value = '  -0.99%  '

$ws.Range("E43").Value = '  -4.74%  '

$ws.Range("D44").Value = '41.65'
$ws.Range("E44").Value = '  -0.25%  '

$ws.Range("E45").Value = '  -1.07%  '

$ws.Range("E46").Value = '  -3.17%  '

$ws.Range("D47").Value = '24.18'
$ws.Range("E47").Value = '  -2.82%  '

$ws.Range("D48").Value = '2.455.34'
$ws.Range("E48").Value = '  +3.45%  '

$ws.Range("B49").Value = 'Cosmos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D49").Value = '6.81'
$ws.Range("E49").Value = '  -2.27%  '

$ws.Range("B50").Value = 'InjectiveProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D50").Value = '23.04'
$ws.Range("E50").Value = '  -1.15%  '

$ws.Range("E51").Value = '  +4.83%  '
